# Update status column for "Apache Commons Imaging" .. "Apache Commons RDF"
# (rows 58-78) from "IN PROGRESS" to "DONE".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C58:C78").Value = "DONE"

# Move the active selection to B58, matching the saved view state.
$ws.Range("B58").Select()
